$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Slit2","Robo2","ECs",3,1,0.2335466666666667,0.70064,0.0710303177517907,0.07103031775179068,3,1,0.2313123333333333,0.693937,0.7569517164947553,0.7569517164947555,0.0540222244088889,0.48620001968,0.05376652094538586,0.05376652094538586),
    @("ECs","Slit2","Robo2","FAPs",3,1,0.2335466666666667,0.70064,0.0710303177517907,0.07103031775179068,1,0.3333333333333333,0.07289866666666667,0.218696,0.2385552472206224,0.2385552472206224,0.01702524060444445,0.15322716544,0.01694465501143779,0.01694465501143779),
    @("ECs","Slit2","Robo2","MuSCs",3,1,0.2335466666666667,0.70064,0.0710303177517907,0.07103031775179068,1,0.3333333333333333,0.001373,0.004119,0.004493036284622232,0.004493036284622232,0.0003206595733333334,0.00288593616,0.0003191417949670422,0.0003191417949670422),
    @("FAPs","Slit2","Robo2","ECs",3,1,2.538883,7.616649,0.772169728639328,0.7721697286393279,3,1,0.2313123333333333,0.693937,0.7569517164947553,0.7569517164947555,0.5872749507903333,5.285474557113,0.5844952015188287,0.5844952015188287),
    @("FAPs","Slit2","Robo2","FAPs",3,1,2.538883,7.616649,0.772169728639328,0.7721697286393279,1,0.3333333333333333,0.07289866666666667,0.218696,0.2385552472206224,0.2385552472206224,0.1850811855226666,1.665730669704,0.1842051405118358,0.1842051405118358),
    @("FAPs","Slit2","Robo2","MuSCs",3,1,2.538883,7.616649,0.772169728639328,0.7721697286393279,1,0.3333333333333333,0.001373,0.004119,0.004493036284622232,0.004493036284622232,0.003485886359,0.031372977231,0.003469386608663403,0.003469386608663403),
    @("MuSCs","Slit2","Robo2","ECs",3,1,0.515556,1.546668,0.1567999536088813,0.1567999536088813,3,1,0.2313123333333333,0.693937,0.7569517164947553,0.7569517164947555,0.119254461324,1.073290151916,0.1186899940305407,0.1186899940305407),
    @("MuSCs","Slit2","Robo2","FAPs",3,1,0.515556,1.546668,0.1567999536088813,0.1567999536088813,1,0.3333333333333333,0.07289866666666667,0.218696,0.2385552472206224,0.2385552472206224,0.037583344992,0.338250104928,0.0374054516973488,0.03740545169734879),
    @("MuSCs","Slit2","Robo2","MuSCs",3,1,0.515556,1.546668,0.1567999536088813,0.1567999536088813,1,0.3333333333333333,0.001373,0.004119,0.004493036284622232,0.004493036284622232,0.0007078583880000001,0.006370725492,0.0007045078809917863,0.0007045078809917862)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}

# Delete the now-obsolete rows 11-13 (old Inflammatory-Mac/MuSCs split collapsed to rows 2-10)
$ws.Range("A11:T13").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp) | Out-Null
